# Append the 2025-02-08 09:03:35 resale-number update as a new row (58)
# on the CityResaleNum sheet, mirroring the existing row layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 58

# Helper: write a value as genuine text (no autodetection into a date/
# number), without leaving behind any number-format/style residue on the
# target cell. We build the literal text in an unused scratch cell via a
# text-producing formula, copy it, and paste-special "values only" into
# the destination - this preserves the destination's existing (default)
# style while still landing a plain text cell.
function Set-TextCell($cellRef, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163) # xlPasteValues
    $scratch.ClearContents()
}

Set-TextCell "A$row" "2025-02-08"
Set-TextCell "B$row" "09:03:35"
Set-TextCell "C$row" "Saturday"
Set-TextCell "D$row" "05"

$ws.Cells.Item($row, 5).Value  = 126202
$ws.Cells.Item($row, 6).Value  = 141818
$ws.Cells.Item($row, 7).Value  = 168031
$ws.Cells.Item($row, 8).Value  = 158275
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 143265
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191572
$ws.Cells.Item($row, 14).Value = 115192
$ws.Cells.Item($row, 15).Value = 44703
$ws.Cells.Item($row, 16).Value = 28318
$ws.Cells.Item($row, 17).Value = 63708
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41113
$ws.Cells.Item($row, 20).Value = -1

$excel.CutCopyMode = $false
